$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "1C B1"
$ws.Range("F3").Value = "1C C1"
$ws.Range("F4").Value = "2C A1"
$ws.Range("F5").Value = "2C B1"
$ws.Range("F6").Value = "2C C1"
$ws.Range("F7").Value = "3C A1"
$ws.Range("F8").Value = "3C B1"
$ws.Range("F9").Value = "3C C1"
$ws.Range("F10").Value = "3C D1"
$ws.Range("F11").Value = "4nC A1"
$ws.Range("F12").Value = "4nC B1"
$ws.Range("F13").Value = "5nC A1"
$ws.Range("F14").Value = "5nC B1"

$ws.Range("F18").Select() | Out-Null
